# Creación de question para la transacción
# Fills the two previously-empty "question" rows (rows 4 and 5) with the
# new "numeroTarjeta / tipoPago / valorPago / moneda / tipoCuenta / numeroCuenta"
# scenario data, then removes the stray leftover cell at M13 (which also
# shrinks the used range back down to A1:T5) and updates the active
# selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "Pago total en dólares" scenario ---------------------------
$ws.Range("O4").Value = "5326666666666666"
$ws.Range("P4").Value = "Pago total en dólares"
$ws.Range("Q4").Value = "0,00"
$ws.Range("R4").Value = "Dólares"
$ws.Range("S4").Value = "Corriente"
$ws.Range("T4").Value = "406-132280-01"

# --- Row 5: "Pago total en pesos" scenario ------------------------------
$ws.Range("O5").Value = "5326666666666666"
$ws.Range("P5").Value = "Pago total en pesos"
$ws.Range("Q5").Value = "0,00"
$ws.Range("R5").Value = "Pesos"
$ws.Range("S5").Value = "Corriente"
$ws.Range("T5").Value = "406-132280-01"

# --- Remove the stray leftover row 13 (only had M13) --------------------
$ws.Rows.Item(13).Delete()

# --- Update the view / selection to match the new layout ---------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 9
$win.ScrollRow = 1
$ws.Range("L13").Select()
